$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Insert a new column before column B ("PubChemCID") for the new
# "ShorthandName" column, shifting all existing data one column to the right.
$ws.Columns.Item(2).Insert()

# Header
$ws.Range("B1").Value = "ShorthandName"

# Shorthand names per compound row (rows 2-13)
$ws.Range("B2").Value = "KI/Cu"
$ws.Range("B3").Value = "Fe/Na/Co"
$ws.Range("B4").Value = "Mg/Mn"
$ws.Range("B5").Value = "V-mix"
$ws.Range("B6").Value = "V-mix2"
$ws.Range("B7").Value = "Ca"
$ws.Range("B8").Value = "Zn"
$ws.Range("B9").Value = "Bor"
$ws.Range("B10").Value = "Ino"
$ws.Range("B11").Value = "Paba"
$ws.Range("B12").Value = "Glu/NH"
$ws.Range("B13").Value = "Buffer"

# Reproduce the saved selection/active cell
$ws.Range("B13").Select()
